# Update cryptos list: refresh D (Price) and E (Volume(1h)) columns
# on the active sheet to reflect the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text would otherwise be auto-parsed as a number
# by Excel (single decimal point, e.g. "213.51") are forced to stay as
# plain text, matching the original inline-string storage. (Each cell is
# formatted individually -- a combined multi-area range only reliably
# applies to the first subrange.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

# Price column (D)
$ws.Range("D2").Value = "27.404.64"
$ws.Range("D3").Value = "1.655.14"
$ws.Range("D5").Value = "213.51"
$ws.Range("D6").Value = "0.513"
$ws.Range("D8").Value = "23.96"
$ws.Range("D11").Value = "0.0879"
$ws.Range("D12").Value = "1.888.50"
$ws.Range("D13").Value = "1.654.72"
$ws.Range("D14").Value = "0.574"
$ws.Range("D15").Value = "4.07"
$ws.Range("D16").Value = "65.72"
$ws.Range("D17").Value = "27.421.11"
$ws.Range("D18").Value = "232.04"
$ws.Range("D20").Value = "7.50"
$ws.Range("D23").Value = "9.31"
$ws.Range("D25").Value = "146.90"
$ws.Range("D27").Value = "15.91"
$ws.Range("D28").Value = "1.00"
$ws.Range("D31").Value = "1.20"
$ws.Range("D33").Value = "1.467.76"
$ws.Range("D37").Value = "0.912"
$ws.Range("D38").Value = "0.573"
$ws.Range("D39").Value = "0.0169"
$ws.Range("D42").Value = "5.45"
$ws.Range("D43").Value = "65.30"
$ws.Range("D45").Value = "1.797.60"
$ws.Range("D46").Value = "0.783"

# Volume(1h) column (E)
$ws.Range("E2").Value = "  -2.20%  "
$ws.Range("E3").Value = "  -1.88%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("E17").Value = "  -2.11%  "
$ws.Range("E18").Value = "  -7.53%  "
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").Value = "  -3.10%  "
$ws.Range("E23").Value = "  -2.40%  "
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -1.43%  "
$ws.Range("E27").Value = "  -3.38%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E32").Value = "  -2.24%  "
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("E35").Value = "  -2.68%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  -2.58%  "
$ws.Range("E38").Value = "  -2.79%  "
$ws.Range("E39").Value = "  -1.47%  "
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -6.38%  "
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("E47").Value = "  +1.05%  "
$ws.Range("E48").Value = "  -1.05%  "
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("E50").Value = "  -1.27%  "
$ws.Range("E51").Value = "  -0.35%  "
